$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4, shifting the existing rows (4-21) down to (5-22).
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new weekly price record.
$ws.Cells.Item(4, 1).Value = 10
$ws.Cells.Item(4, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(4, 3).Value = "La Araucanía"
$ws.Cells.Item(4, 4).Value = 44749
$ws.Cells.Item(4, 5).Value = 9
$ws.Cells.Item(4, 6).Value = 100112017
$ws.Cells.Item(4, 7).Value = "Ramas de apio"
$ws.Cells.Item(4, 8).Value = "Sin especificar"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 65
$ws.Cells.Item(4, 11).Value = 6000
$ws.Cells.Item(4, 12).Value = 6000
$ws.Cells.Item(4, 13).Value = 6000
$ws.Cells.Item(4, 14).Value = "$/paquete"
$ws.Cells.Item(4, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(4, 16).Value = 6000
$ws.Cells.Item(4, 17).Value = 1
$ws.Cells.Item(4, 18).Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Cells.Item(4, 4).NumberFormat = $ws.Cells.Item(5, 4).NumberFormat
